$d = $word.ActiveDocument

$replacements = @(
  @("2025-02-02 Sunday", "2025-02-03 Monday"),
  @("75+8=", "34+42="),
  @("57+26=", "32+19="),
  @("26+54=", "94-24="),
  @("29-22=", "26+55="),
  @("87-64=", "72+7="),
  @("21+36=", "1+86="),
  @("68-6=", "0+69="),
  @("60+19=", "58-52="),
  @("42+25=", "30+63="),
  @("19+34=", "97-28="),
  @("95-75=", "51-4="),
  @("0+61=", "55-8="),
  @("3+51=", "28+20="),
  @("5+87=", "11+53="),
  @("98-89=", "65+5="),
  @("1+45=", "94-48="),
  @("24+11=", "23+12="),
  @("33+62=", "63+32="),
  @("7-7=", "94-5="),
  @("2+24=", "34+23="),
  @("33-13=", "7+89="),
  @("40+38=", "43+19="),
  @("78-35=", "50+13="),
  @("91-72=", "80-18="),
  @("93-55=", "77-68="),
  @("88+7=", "26+8="),
  @("73-8=", "19-18="),
  @("54+39=", "21+45="),
  @("10+37=", "59-18="),
  @("29-13=", "18+81="),
  @("43-4=", "15+13="),
  @("70-62=", "80+11="),
  @("71+26=", "44-9="),
  @("98-10=", "26+68="),
  @("76-8=", "93-45="),
  @("61-47=", "38+18="),
  @("38+2=", "59+14="),
  @("89-46=", "41+34="),
  @("74-46=", "82-52="),
  @("33-10=", "96+0="),
  @("52-30=", "19+42="),
  @("99-93=", "56+32="),
  @("64-3=", "36+3="),
  @("69-66=", "42-29="),
  @("23+9=", "90-36="),
  @("97-17=", "27+45="),
  @("47+9=", "95-1="),
  @("69+1=", "2+59="),
  @("3+68=", "26+55="),
  @("30+41=", "32-9="),
  @("74-19=", "42+36="),
  @("50+36=", "80-53="),
  @("26+7=", "44-7="),
  @("23+76=", "91-34="),
  @("85-18=", "92-63="),
  @("93+6=", "10+46="),
  @("91-87=", "32+46="),
  @("18+66=", "6+73="),
  @("72-66=", "37+59="),
  @("62-51=", "42-21="),
  @("62-13=", "18+77="),
  @("37+9=", "71-22="),
  @("83+13=", "89-11="),
  @("16+49=", "79-66="),
  @("30-14=", "86-10="),
  @("19+45=", "17+10="),
  @("82-55=", "9+31="),
  @("5+3=", "84-17="),
  @("73+22=", "76+19="),
  @("90-11=", "68-65="),
  @("50-16=", "10+59="),
  @("0+9=", "71-10="),
  @("25+69=", "44+28="),
  @("3+88=", "69-8="),
  @("30+28=", "61-8="),
  @("41+18=", "84-56="),
  @("25+16=", "75-57="),
  @("28+19=", "0+14="),
  @("1+60=", "21+66="),
  @("51+24=", "94-6="),
  @("21+51=", "31+17="),
  @("78-65=", "74-24="),
  @("84-37=", "14+78="),
  @("14+5=", "38+7="),
  @("34+31=", "25+50="),
  @("45+20=", "83-63="),
  @("8+25=", "17+45="),
  @("87+0=", "50-24="),
  @("85-10=", "97-36="),
  @("2+86=", "99-33="),
  @("70-9=", "69+11="),
  @("87-75=", "8+58="),
  @("69-3=", "11+48="),
  @("17+34=", "73-25="),
  @("28+22=", "0+30="),
  @("27+12=", "71-70="),
  @("86-12=", "5+66="),
  @("34+61=", "78-3="),
  @("69-35=", "70-4="),
  @("5+2=", "68-56=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
